$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion summary text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 8.4 = 34878.98 pesos`n✅ 34878.98 pesos = 8.36 = 961.9 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$wsHoja1.Range("A1").Value = $newText

# --- Sheet "tasas": update rate cells ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 118.99
$wsTasas.Range("O10").Value = 4150.25
$wsTasas.Range("N12").Value = 4170
$wsTasas.Range("O12").Value = 115.001
